$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1645843333333333
$ws.Range("H2").Value = 0.493753
$ws.Range("I2").Value = 0.03485847193389392
$ws.Range("J2").Value = 0.03485847193389392
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.378475333333333
$ws.Range("N2").Value = 4.135426
$ws.Range("O2").Value = 0.05609715574531157
$ws.Range("P2").Value = 0.05609715574531156
$ws.Range("Q2").Value = 0.2268754437531111
$ws.Range("R2").Value = 2.041878993778
$ws.Range("S2").Value = 0.001955461129119219
$ws.Range("T2").Value = 0.001955461129119219

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1645843333333333
$ws.Range("H3").Value = 0.493753
$ws.Range("I3").Value = 0.03485847193389392
$ws.Range("J3").Value = 0.03485847193389392
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.16176133333333
$ws.Range("N3").Value = 39.485284
$ws.Range("O3").Value = 0.5356188518899525
$ws.Range("P3").Value = 0.5356188518899525
$ws.Range("Q3").Value = 2.166219714539111
$ws.Range("R3").Value = 19.495977430852
$ws.Range("S3").Value = 0.01867085471587039
$ws.Range("T3").Value = 0.01867085471587039

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1645843333333333
$ws.Range("H4").Value = 0.493753
$ws.Range("I4").Value = 0.03485847193389392
$ws.Range("J4").Value = 0.03485847193389392
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.829094
$ws.Range("N4").Value = 2.487282
$ws.Range("O4").Value = 0.03374004171190829
$ws.Range("P4").Value = 0.03374004171190828
$ws.Range("Q4").Value = 0.1364558832606667
$ws.Range("R4").Value = 1.228102949346
$ws.Range("S4").Value = 0.001176126297062965
$ws.Range("T4").Value = 0.001176126297062965

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1645843333333333
$ws.Range("H5").Value = 0.493753
$ws.Range("I5").Value = 0.03485847193389392
$ws.Range("J5").Value = 0.03485847193389392
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.203668
$ws.Range("N5").Value = 27.611004
$ws.Range("O5").Value = 0.3745439506528278
$ws.Range("P5").Value = 0.3745439506528276
$ws.Range("Q5").Value = 1.514779562001333
$ws.Range("R5").Value = 13.633016058012
$ws.Range("S5").Value = 0.01305602979184134
$ws.Range("T5").Value = 0.01305602979184134

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.368329
$ws.Range("H6").Value = 10.104987
$ws.Range("I6").Value = 0.7134020567608963
$ws.Range("J6").Value = 0.7134020567608964
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.378475333333333
$ws.Range("N6").Value = 4.135426
$ws.Range("O6").Value = 0.05609715574531157
$ws.Range("P6").Value = 0.05609715574531156
$ws.Range("Q6").Value = 4.643158441051333
$ws.Range("R6").Value = 41.788425969462
$ws.Range("S6").Value = 0.0400198262871416
$ws.Range("T6").Value = 0.0400198262871416

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.368329
$ws.Range("H7").Value = 10.104987
$ws.Range("I7").Value = 0.7134020567608963
$ws.Range("J7").Value = 0.7134020567608964
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.16176133333333
$ws.Range("N7").Value = 39.485284
$ws.Range("O7").Value = 0.5356188518899525
$ws.Range("P7").Value = 0.5356188518899525
$ws.Range("Q7").Value = 44.33314239014533
$ws.Range("R7").Value = 398.998281511308
$ws.Range("S7").Value = 0.382111590578202
$ws.Range("T7").Value = 0.382111590578202

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.368329
$ws.Range("H8").Value = 10.104987
$ws.Range("I8").Value = 0.7134020567608963
$ws.Range("J8").Value = 0.7134020567608964
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.829094
$ws.Range("N8").Value = 2.487282
$ws.Range("O8").Value = 0.03374004171190829
$ws.Range("P8").Value = 0.03374004171190828
$ws.Range("Q8").Value = 2.792661363926
$ws.Range("R8").Value = 25.133952275334
$ws.Range("S8").Value = 0.0240702151524738
$ws.Range("T8").Value = 0.0240702151524738

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.368329
$ws.Range("H9").Value = 10.104987
$ws.Range("I9").Value = 0.7134020567608963
$ws.Range("J9").Value = 0.7134020567608964
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.203668
$ws.Range("N9").Value = 27.611004
$ws.Range("O9").Value = 0.3745439506528278
$ws.Range("P9").Value = 0.3745439506528276
$ws.Range("Q9").Value = 31.000981830772
$ws.Range("R9").Value = 279.008836476948
$ws.Range("S9").Value = 0.267200424743079
$ws.Range("T9").Value = 0.2672004247430789

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.188588333333333
$ws.Range("H10").Value = 3.565765
$ws.Range("I10").Value = 0.2517394713052097
$ws.Range("J10").Value = 0.2517394713052098
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.378475333333333
$ws.Range("N10").Value = 4.135426
$ws.Range("O10").Value = 0.05609715574531157
$ws.Range("P10").Value = 0.05609715574531156
$ws.Range("Q10").Value = 1.638439698987778
$ws.Range("R10").Value = 14.74595729089
$ws.Range("S10").Value = 0.01412186832905074
$ws.Range("T10").Value = 0.01412186832905074

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.188588333333333
$ws.Range("H11").Value = 3.565765
$ws.Range("I11").Value = 0.2517394713052097
$ws.Range("J11").Value = 0.2517394713052098
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 13.16176133333333
$ws.Range("N11").Value = 39.485284
$ws.Range("O11").Value = 0.5356188518899525
$ws.Range("P11").Value = 0.5356188518899525
$ws.Range("Q11").Value = 15.64391596691778
$ws.Range("R11").Value = 140.79524370226
$ws.Range("S11").Value = 0.1348364065958801
$ws.Range("T11").Value = 0.1348364065958801

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.188588333333333
$ws.Range("H12").Value = 3.565765
$ws.Range("I12").Value = 0.2517394713052097
$ws.Range("J12").Value = 0.2517394713052098
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.829094
$ws.Range("N12").Value = 2.487282
$ws.Range("O12").Value = 0.03374004171190829
$ws.Range("P12").Value = 0.03374004171190828
$ws.Range("Q12").Value = 0.9854514556366667
$ws.Range("R12").Value = 8.869063100729999
$ws.Range("S12").Value = 0.008493700262371517
$ws.Range("T12").Value = 0.008493700262371517

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.188588333333333
$ws.Range("H13").Value = 3.565765
$ws.Range("I13").Value = 0.2517394713052097
$ws.Range("J13").Value = 0.2517394713052098
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.203668
$ws.Range("N13").Value = 27.611004
$ws.Range("O13").Value = 0.3745439506528278
$ws.Range("P13").Value = 0.5356188518899525
$ws.Range("Q13").Value = 10.93937240867333
$ws.Range("R13").Value = 98.45435167806001
$ws.Range("S13").Value = 0.09428749611790742
$ws.Range("T13").Value = 0.09428749611790742
